# Apply the commit's changes to the workbook:
#  1. Remove the (empty) B4 / B5 cells on the "ODI Batting" sheet.
#  2. Add a new "ODI Batting Extra" worksheet (after "ODI Bowling") with the
#     MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#     MAN_OF_MATCH table.

$wb = $excel.ActiveWorkbook

# --- 1. "ODI Batting" sheet: clear the stray empty B4 / B5 cells ---------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B4").ClearContents()
$batting.Range("B5").ClearContents()

# --- 2. New "ODI Batting Extra" sheet, appended after "ODI Bowling" ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row: re-use the bold / bordered / centered header style already
# used by the other sheets (e.g. "ODI Batting" row 1).
$batting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

# Touch every cell in the A1:F10 block up front (matches the source data
# dump, which materialises a blank placeholder cell for every row/column
# even where there is no value). MATCH_CODE (A), NUM_4 (C), NUM_6 (D) and
# PERCENT_RUNS_OF_TOTAL (E) are all stored as text, like the analogous
# MATCH_CODE column on the other sheets - force text formatting up front so
# the values aren't re-interpreted as numbers/percentages. BATTING_POSITION
# (B) and MAN_OF_MATCH (F) keep the default/general format.
$extra.Range("A2:A10").NumberFormat = "@"
$extra.Range("B2:B10").NumberFormat = "General"
$extra.Range("C2:E10").NumberFormat = "@"
$extra.Range("F2:F10").NumberFormat = "General"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

$extra.Range("A2").Value = "3526"
$extra.Range("B2").Value = 9
$extra.Range("C2").Value = "0"
$extra.Range("D2").Value = "0"
$extra.Range("E2").Value = "2.03%"
$extra.Range("F2").Value = "NO"

$extra.Range("A3").Value = "3528"
$extra.Range("B3").Value = 9
$extra.Range("C3").Value = "1"
$extra.Range("D3").Value = "1"
$extra.Range("E3").Value = "5.22%"
$extra.Range("F3").Value = "NO"

$extra.Range("A4").Value = "3605"
$extra.Range("B4").Value = 7
$extra.Range("F4").Value = "YES"

$extra.Range("A5").Value = "3610"
$extra.Range("B5").Value = 6
$extra.Range("C5").Value = "0"
$extra.Range("D5").Value = "0"
$extra.Range("E5").Value = "1.40%"
$extra.Range("F5").Value = "NO"

$extra.Range("A6").Value = "4184"
$extra.Range("F6").Value = "NO"

$extra.Range("A7").Value = "4185"

$extra.Range("A8").Value = "4563"

$extra.Range("A9").Value = "4566"

$extra.Range("A10").Value = "4568"

[void]$extra.Range("A1").Select()
